# 운동기록.xlsx edit script
# - rename sheet1, add SQL sheet with its query string
# - mark sheet1 "fit to page"
# - append a new data row (7/1/2020) to sheet1

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Page setup: fit to page
$ws1.PageSetup.FitToPagesWide = 1
$ws1.PageSetup.FitToPagesTall = 1

# Append new row of exercise data (copy formats from the row above, then set values)
$ws1.Range("A27:L27").Copy()
$ws1.Range("A28:L28").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws1.Range("A28").Value = 44013
$ws1.Range("B28").Value = 97.5
$ws1.Range("C28").Value = 104.5
$ws1.Range("D28").Value = 0.93
$ws1.Range("E28").Value = "SAME"
$ws1.Range("F28").Value = 0
$ws1.Range("G28").Value = 82.6
$ws1.Range("H28").Value = 0.58
$ws1.Range("I28").Value = "SAME"
$ws1.Range("J28").Value = 0
$ws1.Range("K28").Value = 28.9
$ws1.Range("L28").Value = "OVERWEIGHT"

# Rename sheet1 to reflect the source-table identity used by the SQL sheet
$ws1.Name = 'SEUNGJAE_HAN."운동기록"'

# Add a second sheet holding the SQL query used to build this report
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "SQL"

$sql = 'select DATETIME DATETIME, WAIST WAIST, HIP HIP, WHR WHR, WHR_IMPROVEMENT WHR_IMPROVEMENT, WHR_CHANGE WHR_CHANGE, WEIGHT WEIGHT, WHTR WHTR, WHTR_IMPROVEMENT WHTR_IMPROVEMENT, WHTR_CHANGE WHTR_CHANGE, BMI BMI, OBESITY OBESITY from (select * from "SEUNGJAE_HAN"."운동기록")'
$ws2.Range("A2").Value = $sql

$ws1.Activate()
